$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sending cluster (column A) / Target cluster (column D) relabeling ---
# New shared string "ECs" cluster appears; rows 2-3 sending cluster becomes
# "ECs" (was "FAPs"), rows 4-5 sending cluster becomes "FAPs" (was "MuSCs").
# Target cluster (column D) keeps its semantic FAPs/MuSCs meaning.
$ws.Range("A2").Value = "ECs"
$ws.Range("A3").Value = "ECs"
$ws.Range("A4").Value = "FAPs"
$ws.Range("A5").Value = "FAPs"

$ws.Range("D2").Value = "FAPs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("D5").Value = "MuSCs"

# --- Row 2 numeric updates ---
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01364166666666667
$ws.Range("H2").Value = 0.040925
$ws.Range("I2").Value = 0.005793924852122192
$ws.Range("J2").Value = 0.005793924852122191
$ws.Range("M2").Value = 0.7050613333333334
$ws.Range("O2").Value = 0.6029366303164088
$ws.Range("P2").Value = 0.6949112282957692
$ws.Range("Q2").Value = 0.009618211688888892
$ws.Range("R2").Value = 0.08656390520000001
$ws.Range("S2").Value = 0.003493369526645051
$ws.Range("T2").Value = 0.004026263435641615

# --- Row 3 numeric updates ---
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01364166666666667
$ws.Range("H3").Value = 0.040925
$ws.Range("I3").Value = 0.005793924852122192
$ws.Range("J3").Value = 0.005793924852122191
$ws.Range("M3").Value = 0.4643175
$ws.Range("N3").Value = 0.928635
$ws.Range("O3").Value = 0.3970633696835912
$ws.Range("P3").Value = 0.3050887717042308
$ws.Range("Q3").Value = 0.006334064562500001
$ws.Range("R3").Value = 0.038004387375
$ws.Range("S3").Value = 0.002300555325477141
$ws.Range("T3").Value = 0.001767661416480576

# --- Row 4 numeric updates ---
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.340836
$ws.Range("H4").Value = 7.022508
$ws.Range("I4").Value = 0.9942060751478778
$ws.Range("J4").Value = 0.9942060751478778
$ws.Range("M4").Value = 0.7050613333333334
$ws.Range("O4").Value = 0.6029366303164088
$ws.Range("P4").Value = 0.6949112282957692
$ws.Range("Q4").Value = 1.650432951274667
$ws.Range("R4").Value = 14.853896561472
$ws.Range("S4").Value = 0.5994432607897637
$ws.Range("T4").Value = 0.6908849648601276

# --- Row 5 numeric updates ---
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.340836
$ws.Range("H5").Value = 7.022508
$ws.Range("I5").Value = 0.9942060751478778
$ws.Range("J5").Value = 0.9942060751478778
$ws.Range("M5").Value = 0.4643175
$ws.Range("N5").Value = 0.928635
$ws.Range("O5").Value = 0.3970633696835912
$ws.Range("P5").Value = 0.3050887717042308
$ws.Range("Q5").Value = 1.08689111943
$ws.Range("R5").Value = 6.52134671658
$ws.Range("S5").Value = 0.3947628143581141
$ws.Range("T5").Value = 0.3033211102877502
